# Apply cell text updates from the crypto price refresh diff.
# Each cell is forced to remain a text string (matching the source
# workbook's inlineStr cells) by briefly switching to a Text number
# format before the assignment, then clearing formatting again so no
# stray style index is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@('D2', '27.520.80')
    ,@('E2', '  +1.09%  ')
    ,@('D3', '1.767.83')
    ,@('E3', '  -0.96%  ')
    ,@('D4', '1.010')
    ,@('E4', '  +0.66%  ')
    ,@('D5', '337.71')
    ,@('E5', '  +0.92%  ')
    ,@('D6', '1.006')
    ,@('E6', '  +0.51%  ')
    ,@('D7', '0.3859')
    ,@('E7', '  +1.99%  ')
    ,@('D8', '0.3424')
    ,@('E8', '  -0.33%  ')
    ,@('D9', '46.74')
    ,@('E9', '  -3.39%  ')
    ,@('D10', '1.148')
    ,@('D11', '0.07449')
    ,@('E11', '  -0.70%  ')
    ,@('D12', '1.008')
    ,@('E12', '  +0.62%  ')
    ,@('D13', '22.46')
    ,@('E13', '  +3.13%  ')
    ,@('E14', '  -1.62%  ')
    ,@('D15', '1.768.37')
    ,@('E15', '  -1.22%  ')
    ,@('D16', '7.080')
    ,@('E16', '  -0.38%  ')
    ,@('D17', '0.00001077')
    ,@('E17', '  -1.99%  ')
    ,@('D18', '0.06692')
    ,@('E18', '  +0.32%  ')
    ,@('D19', '82.27')
    ,@('E19', '  -1.85%  ')
    ,@('D20', '1.006')
    ,@('E20', '  +0.53%  ')
    ,@('D21', '17.46')
    ,@('E21', '  +0.55%  ')
    ,@('D22', '6.463')
    ,@('E22', '  -2.44%  ')
    ,@('D23', '27.544.53')
    ,@('E23', '  +1.17%  ')
    ,@('D24', '12.20')
    ,@('E24', '  -1.60%  ')
    ,@('D25', '2.379')
    ,@('E25', '  -1.62%  ')
    ,@('B26', 'ImmutableX')
    ,@('C26', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx')
    ,@('D26', '1.447')
    ,@('E26', '  -3.68%  ')
    ,@('B27', 'EthereumClassic')
    ,@('C27', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc')
    ,@('D27', '20.82')
    ,@('E27', '  -2.43%  ')
    ,@('D28', '2.448')
    ,@('E28', '  -4.00%  ')
    ,@('D29', '153.38')
    ,@('E29', '  -0.33%  ')
    ,@('B30', 'BitcoinCash')
    ,@('C30', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch')
    ,@('D30', '134.77')
    ,@('E30', '  +0.47%  ')
    ,@('B31', 'WrappedliquidstakedEther2.0')
    ,@('C31', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth')
    ,@('D31', '1.967.77')
    ,@('E31', '  -1.26%  ')
    ,@('D32', '6.159')
    ,@('E32', '  +0.72%  ')
    ,@('D33', '3.968')
    ,@('E33', '  -1.34%  ')
    ,@('D34', '0.08849')
    ,@('E34', '  +1.69%  ')
    ,@('D35', '12.70')
    ,@('E35', '  -4.44%  ')
    ,@('D36', '0.02440')
    ,@('E36', '  +4.34%  ')
    ,@('D37', '5.404')
    ,@('E37', '  -1.08%  ')
    ,@('D38', '0.6805')
    ,@('E38', '  -2.22%  ')
    ,@('D39', '0.06352')
    ,@('E39', '  +0.26%  ')
    ,@('D40', '0.2203')
    ,@('E40', '  -0.11%  ')
    ,@('D41', '1.547')
    ,@('E41', '  -6.82%  ')
    ,@('D42', '1.246')
    ,@('E42', '  +0.26%  ')
    ,@('D43', '8.395')
    ,@('E43', '  -4.90%  ')
    ,@('E44', '  -1.20%  ')
    ,@('D45', '1.006')
    ,@('E45', '  +0.52%  ')
    ,@('D46', '0.6266')
    ,@('E46', '  -3.86%  ')
    ,@('D47', '3.846')
    ,@('E47', '  -0.05%  ')
    ,@('D48', '132.14')
    ,@('D49', '2.114')
    ,@('E49', '  -1.62%  ')
    ,@('D50', '0.07410')
    ,@('E50', '  +3.76%  ')
    ,@('E51', '  +2.48%  ')
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.ClearFormats()
}
